# "Atualização dos requisitos" - update requirements wording, widen the
# Descrição column, move the selection, and shrink the print scale.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Requirement text updates (column D / "Descrição") --------------------

# RF02: clarify that the park/trail list lives under the PARQUES tab.
$ws.Range("D3").Value = "Exibir lista de parques, trilhas e pontos turísticos com imagens e detalhes na aba PARQUES (PARNASO, PETP, Montanhas de Teresópolis)."

# RNF02: tie the dark-green palette to the local flora.
$ws.Range("D9").Value = "Manter interface clara, intuitiva, legível e coerente com paleta de verde escuro, fazendo uma alusão a flora de Teresópolis."

# RNF03: replaced by a new requirement about the footer logo.
$ws.Range("D10").Value = "Apresentar a logo do Terê Verde também ao rodapé da página. "

# --- Layout tweaks ----------------------------------------------------------

# Widen column D (Descrição) to fit the longer text.
$ws.Columns.Item(4).ColumnWidth = 146.02213541666666

# Move the active selection to E23.
$ws.Range("E23").Select() | Out-Null

# Shrink the print scale from 58% to 55%.
$ws.PageSetup.Zoom = 55
